$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    # Force the cell to be stored as text (matching the workbook's existing
    # inline-string / text cells), even when the new value looks numeric.
    # The leading apostrophe is Excel's quote-prefix for literal text, and
    # resetting the Style back to "Normal" afterwards avoids leaving a
    # Text (@) number-format style attached to the cell.
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).Style = "Normal"
}

# Row 41/42 special-case: the two rows are swapped (Aptos moves up to row
# 41, FraxShare moves down to row 42) and both get refreshed price /
# volume figures.
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell "D41" "10.49"
$ws.Range("E41").Value = "  -1.12%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell "D42" "7.661"
$ws.Range("E42").Value = "  +4.67%  "

# Price (column D) and Volume(1h) (column E) refreshes for every other row.
$updates = @{
    2  = @{ D = "20.593.51";   E = "  +0.52%  " }
    3  = @{ D = "1.480.40";    E = "  +0.80%  " }
    4  = @{ D = "1.010";       E = "  +0.07%  " }
    5  = @{ D = "0.9725";      E = "  +1.70%  " }
    6  = @{ D = "279.32";      E = "  -0.73%  " }
    7  = @{ D = "0.3664";      E = "  -0.96%  " }
    8  = @{ D = "0.3081";      E = "  -2.98%  " }
    9  = @{ D = "40.02";       E = "  -4.27%  " }
    10 = @{ D = "1.066";       E = "  +0.92%  " }
    11 = @{ D = "0.06678";     E = "  +0.06%  " }
    12 = @{ E = "  -0.06%  " }
    13 = @{ D = "5.528";       E = "  -1.43%  " }
    14 = @{ D = "18.10";       E = "  -0.50%  " }
    15 = @{ D = "6.220";       E = "  -0.52%  " }
    16 = @{ D = "0.9735";      E = "  +1.74%  " }
    17 = @{ D = "0.00001030";  E = "  -0.55%  " }
    18 = @{ D = "1.476.82";    E = "  +0.23%  " }
    19 = @{ D = "0.05944";     E = "  +4.18%  " }
    20 = @{ D = "69.76";       E = "  -3.35%  " }
    21 = @{ D = "5.502";       E = "  -2.83%  " }
    22 = @{ D = "14.54";       E = "  -1.12%  " }
    23 = @{ D = "11.06";       E = "  -1.43%  " }
    24 = @{ D = "2.257";       E = "  -0.43%  " }
    25 = @{ D = "20.638.76";   E = "  -0.21%  " }
    26 = @{ D = "142.17";      E = "  +3.04%  " }
    27 = @{ D = "2.135";       E = "  -6.72%  " }
    28 = @{ D = "17.29";       E = "  -1.41%  " }
    29 = @{ D = "1.637.14";    E = "  -0.04%  " }
    30 = @{ D = "114.28";      E = "  +0.33%  " }
    31 = @{ D = "3.932";       E = "  -0.48%  " }
    32 = @{ D = "5.030";       E = "  -5.46%  " }
    33 = @{ D = "0.8232";      E = "  -1.17%  " }
    34 = @{ E = "  +2.08%  " }
    35 = @{ D = "1.537";       E = "  -5.12%  " }
    36 = @{ D = "1.205";       E = "  +7.48%  " }
    37 = @{ D = "0.05783";     E = "  -4.09%  " }
    38 = @{ D = "4.726";       E = "  -3.73%  " }
    39 = @{ D = "0.9733";      E = "  +0.37%  " }
    40 = @{ D = "0.02047";     E = "  -1.19%  " }
    43 = @{ D = "0.1885";      E = "  -0.06%  " }
    44 = @{ D = "0.5306";      E = "  -1.86%  " }
    45 = @{ D = "3.536";       E = "  -1.55%  " }
    46 = @{ D = "12.26";       E = "  -1.68%  " }
    47 = @{ D = "118.52";      E = "  -2.99%  " }
    48 = @{ D = "0.5200";      E = "  -2.07%  " }
    49 = @{ D = "1.809";       E = "  -1.09%  " }
    50 = @{ E = "  +0.75%  " }
    51 = @{ D = "0.9953";      E = "  +0.10%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        Set-TextCell "D$row" $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        $ws.Range("E$row").Value = $vals["E"]
    }
}
